# Updated symbol list on Sat Dec 24 19:33:32 UTC 2022 with GitHub Actions
#
# This script reproduces price/volume refresh edits plus a 3-row
# reshuffle (rows 41-43) in the crypto price table on Sheet1.
#
# Numeric-looking text values are written with a leading apostrophe so
# that Excel keeps them as text (preserving exact formatting such as
# leading/trailing zeros) instead of silently coercing them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) refreshes -----------------------------------
$ws.Range("D2").Value  = "'244.77"
$ws.Range("D4").Value  = "'5.395"
$ws.Range("D5").Value  = "'0.06039"
$ws.Range("D6").Value  = "'3.391"
$ws.Range("D7").Value  = "'0.8141"
$ws.Range("D8").Value  = "'0.9336"
$ws.Range("D10").Value = "'0.07497"
$ws.Range("D11").Value = "'0.03519"
$ws.Range("D12").Value = "'0.03062"
$ws.Range("D13").Value = "'0.09437"
$ws.Range("D14").Value = "'4.015"
$ws.Range("D18").Value = "'0.005641"
$ws.Range("D20").Value = "'0.0009921"
$ws.Range("D21").Value = "'3.669"
$ws.Range("D22").Value = "'6.422"
$ws.Range("D23").Value = "'2.180"
$ws.Range("D26").Value = "'0.00007002"
$ws.Range("D40").Value = "'0.03999"

# --- Rows 41-43 reshuffle (BKEXToken/CEJI/KickToken rotate) -------
# Row 41 becomes KickToken (was BKEXToken)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006402"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 becomes BKEXToken (was CEJI)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1079"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 becomes CEJI (was KickToken)
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002901"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining price (column D) refreshes --------------------------
$ws.Range("D44").Value = "'0.005920"
$ws.Range("D45").Value = "'0.00005235"
$ws.Range("D47").Value = "'1.0000"
$ws.Range("D48").Value = "'0.002359"
$ws.Range("D49").Value = "'0.00002101"
